$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$used = $ws.UsedRange
$lastRow = $used.Rows.Count

$changedCount = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null -and $val.StartsWith("System")) {
        $parts = $val -split ", "
        if ($parts.Count -gt 1) {
            $reversed = $parts[($parts.Count - 1)..0]
            $newVal = $reversed -join ", "
            $cell.Value2 = $newVal
            $changedCount = $changedCount + 1
        }
    }
}

Write-Host ("Recorded By column normalized; cells changed: " + $changedCount)
